$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Inflammatory-Mac" target-cluster row (old row 4). This shifts the
# "MuSCs" row up to row 4 and the "Resolving-Mac" row up to row 5; the
# now-unused "Inflammatory-Mac" shared string is dropped automatically.
$ws.Rows(4).Delete()

# Refresh the receptor / edge-weight columns (K:T) with the updated
# TPM-derived numbers from the new script run, for every remaining data row.

# Row 2 (Target cluster: ECs)
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.7195943333333332
$ws.Range("N2").Value = 2.158783
$ws.Range("O2").Value = 0.03204779321415739
$ws.Range("P2").Value = 0.03468223907394029
$ws.Range("Q2").Value = 0.1350537043447778
$ws.Range("R2").Value = 1.215483339103
$ws.Range("S2").Value = 0.03204779321415739
$ws.Range("T2").Value = 0.03468223907394029

# Row 3 (Target cluster: FAPs)
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 16.53477466666667
$ws.Range("N3").Value = 49.60432400000001
$ws.Range("O3").Value = 0.7363913455312854
$ws.Range("P3").Value = 0.7969254084681946
$ws.Range("Q3").Value = 3.103252021031556
$ws.Range("R3").Value = 27.929268189284
$ws.Range("S3").Value = 0.7363913455312854
$ws.Range("T3").Value = 0.7969254084681946

# Row 4 (Target cluster: MuSCs, previously row 5)
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.116736
$ws.Range("N4").Value = 10.233472
$ws.Range("O4").Value = 0.2278785277530463
$ws.Range("P4").Value = 0.1644073176694804
$ws.Range("Q4").Value = 0.9603107180586666
$ws.Range("R4").Value = 5.761864308351999
$ws.Range("S4").Value = 0.2278785277530463
$ws.Range("T4").Value = 0.1644073176694804

# Row 5 (Target cluster: Resolving-Mac, previously row 6)
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.08268233333333334
$ws.Range("N5").Value = 0.248047
$ws.Range("O5").Value = 0.003682333501510851
$ws.Range("P5").Value = 0.003985034788384784
$ws.Range("Q5").Value = 0.01551784788077778
$ws.Range("R5").Value = 0.139660630927
$ws.Range("S5").Value = 0.003682333501510851
$ws.Range("T5").Value = 0.003985034788384784
